$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.641.17"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "1.788.19"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'223.01"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("D6").Value = "'0.556"
$ws.Range("E6").Value = "  -0.91%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'32.65"
$ws.Range("E8").Value = "  +6.94%  "
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").Value = "'0.0679"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Value = "'0.0937"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "2.044.17"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "'11.17"
$ws.Range("E13").Value = "  +10.82%  "
$ws.Range("D14").Value = "1.782.33"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "34.610.45"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Value = "'68.48"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'253.22"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "0.0₃0773"
$ws.Range("E20").Value = "  +4.66%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "'4.21"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'158.75"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'16.34"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = "  +1.90%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "'3.57"
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "1.440.77"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  +4.38%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").Value = "'5.95"
$ws.Range("E45").Value = "  +3.90%  "
$ws.Range("D46").Value = "'1.04"
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("D47").Value = "1.944.38"
$ws.Range("D48").Value = "'105.02"
$ws.Range("E48").Value = "  +7.61%  "
$ws.Range("D49").Value = "'12.03"
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'49.55"
$ws.Range("E51").Value = "  -2.41%  "
